{"js": "// Load all paragraphs in the document body so we can locate the anchor\n// paragraphs by their text content.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Edit 1 -------------------------------------------------------------\n// Insert a new bold paragraph right after the \"Since his arrest in 2012...\"\n// paragraph (and therefore right before the \"MLT was a part of an\n// anti-Israel anarchist group...\" paragraph), describing the unverified\n// Lizard Squad connection found in a dox dump.\nconst arrestParaText = \"Since his arrest in 2012\";\nlet arrestPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(arrestParaText) !== -1) {\n    arrestPara = p;\n    break;\n  }\n}\nif (!arrestPara) {\n  throw new Error(\"Could not find the 'Since his arrest in 2012...' paragraph\");\n}\n\nconst lsConnectionText =\n  \"There is an unverified and potentially false connection via a dox dump \" +\n  \"that states that MLT was a new member of Lizard Squad at one time. It is \" +\n  \"unknown if he played any large part in LS or if he was part of it at \" +\n  \"all, as there is very little information in the dump about MLT.\";\n\nconst lsConnectionPara = arrestPara.insertParagraph(\n  lsConnectionText,\n  Word.InsertLocation.after\n);\nlsConnectionPara.font.bold = true;\n\n// --- Edit 2 -------------------------------------------------------------\n// Insert a new \"Lizard Squad\" bullet in the Connections list, right after\n// the \"Mujihadeen Hacker Group\" bullet.\nconst mujText = \"Mujihadeen Hacker Group\";\nlet mujPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(mujText) !== -1) {\n    mujPara = p;\n    break;\n  }\n}\nif (!mujPara) {\n  throw new Error(\"Could not find the 'Mujihadeen Hacker Group' paragraph\");\n}\n\nconst lizardSquadPara = mujPara.insertParagraph(\n  \"Lizard Squad\",\n  Word.InsertLocation.after\n);\nlizardSquadPara.font.bold = true;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------------\n# Insert a new bold paragraph right after the \"Since his arrest in 2012...\"\n# paragraph (and therefore right before the \"MLT was a part of an\n# anti-Israel anarchist group...\" paragraph), describing the unverified\n# Lizard Squad connection found in a dox dump.\n$findRange = $d.Content\n$findRange.Find.Text = \"Since his arrest in 2012\"\n$found = $findRange.Find.Execute()\nif (-not $found) {\n    throw \"Could not find the 'Since his arrest in 2012...' paragraph\"\n}\n$arrestPara = $findRange.Paragraphs(1)\n$arrestPara.Range.InsertParagraphAfter()\n$lsConnectionPara = $arrestPara.Next()\n$lsConnectionPara.Range.Text = \"There is an unverified and potentially false connection via a dox dump that states that MLT was a new member of Lizard Squad at one time. It is unknown if he played any large part in LS or if he was part of it at all, as there is very little information in the dump about MLT.\"\n$lsConnectionPara.Range.Font.Bold = $true\n\n# --- Edit 2 ---------------------------------------------------------------\n# Insert a new \"Lizard Squad\" bullet in the Connections list, right after\n# the \"Mujihadeen Hacker Group\" bullet.\n$findRange2 = $d.Content\n$findRange2.Find.Text = \"Mujihadeen Hacker Group\"\n$found2 = $findRange2.Find.Execute()\nif (-not $found2) {\n    throw \"Could not find the 'Mujihadeen Hacker Group' paragraph\"\n}\n$mujPara = $findRange2.Paragraphs(1)\n$mujPara.Range.InsertParagraphAfter()\n$lizardSquadPara = $mujPara.Next()\n$lizardSquadPara.Range.Text = \"Lizard Squad\"\n$lizardSquadPara.Range.Font.Bold = $true\n"}
